$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Update the "no claim found" body message (B6)
$ws.Range("B6").Value = "`nBegin body message:`nThere is no claim found`n"

# Update the footer / signature text (B7) - drop the phone number line
$ws.Range("B7").Value = "Best Regards, `nNhut Dang `nAS White Global `nAustralia | Vietnam `nREE Tower, Level 7, 9 Doan Van Bo, Ward 9, Dictrict 4, HCMC, Vietnam `n"

# The shorter B6 text now wraps to fewer lines, so the row shrinks
$ws.Rows.Item(6).RowHeight = 100.8

# Update the selection on Sheet2 to match the saved view state
$ws.Range("E6").Select()
